$d = $word.ActiveDocument

$replacements = @(
    @("67×53=", "44×36="),
    @("40×72=", "42×72="),
    @("60×53=", "23×20="),
    @("91×53=", "28×83="),
    @("95×98=", "79×68="),
    @("84×36=", "88×18="),
    @("65×85=", "33×88="),
    @("40×55=", "54×97="),
    @("75×43=", "14×21="),
    @("79×77=", "50×91="),
    @("25×85=", "97×25="),
    @("89×84=", "76×98="),
    @("12×85=", "65×46="),
    @("18×88=", "89×26="),
    @("17×91=", "74×72="),
    @("24×36=", "65×54="),
    @("94×25=", "16×81="),
    @("62×71=", "98×79="),
    @("39×26=", "81×67="),
    @("17×41=", "81×94="),
    @("24×55=", "11×32="),
    @("61×11=", "21×26="),
    @("34×81=", "99×68="),
    @("64×69=", "85×88="),
    @("78×50=", "64×85=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
